# Assignments.xlsx edit:
# The single "Practicals" row (Java Programming, "Ongoing (4 total)", 30%) is
# replaced by four individual practical rows (Practical 2/3/5/7), each with a
# real due date and a 7.5% weight. Everything below the Java Programming
# block shifts down by three rows to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert 3 new rows right below the existing "Practicals" row
# (old row 10). This pushes "Programming Project" (old row 11) and
# "In Class Test" (old row 12), plus everything further down, down by 3 rows.
$ws.Rows("11:13").Insert()

# --- Grab a ready-made date format (from the still-unshifted SRS date cell,
# now at C18) and stamp it (format only) onto the four new date cells so we
# reuse the workbook's existing date style instead of minting a new one.
$ws.Cells.Item(18, 3).Copy()
$ws.Range("C10:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 10: was "Practicals" / Java Programming / Ongoing (4 total) / 0.3
#     becomes "Practical 2" / Java Programming / 11 Feb 2021 / 0.075
$ws.Cells.Item(10, 1).Value2 = "Practical 2"
$ws.Cells.Item(10, 2).Value2 = "Java Programming"
$ws.Cells.Item(10, 3).Value2 = 44238
$ws.Cells.Item(10, 4).Value2 = 0.075

# --- Row 11 (new): "Practical 3" / Java Programming / 18 Feb 2022 / 0.075
$ws.Cells.Item(11, 1).Value2 = "Practical 3"
$ws.Cells.Item(11, 2).Value2 = "Java Programming"
$ws.Cells.Item(11, 3).Value2 = 44610
$ws.Cells.Item(11, 4).Value2 = 0.075

# --- Row 12 (new): "Practical 5" / Java Programming / 4 Mar 2022 / 0.075
$ws.Cells.Item(12, 1).Value2 = "Practical 5"
$ws.Cells.Item(12, 2).Value2 = "Java Programming"
$ws.Cells.Item(12, 3).Value2 = 44624
$ws.Cells.Item(12, 4).Value2 = 0.075

# --- Row 13 (new): "Practical 7" / Java Programming / 18 Mar 2022 / 0.075
$ws.Cells.Item(13, 1).Value2 = "Practical 7"
$ws.Cells.Item(13, 2).Value2 = "Java Programming"
$ws.Cells.Item(13, 3).Value2 = 44638
$ws.Cells.Item(13, 4).Value2 = 0.075

# --- Update the selection to land on C14, matching the saved workbook state.
$ws.Range("C14").Select() | Out-Null
